$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to Text format so numeric-looking strings
# like "1.003" are stored as text, matching the source data's inlineStr type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.766.68"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "1.743.16"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "333.22"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "0.3886"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("D8").Value = "0.3369"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").Value = "45.27"
$ws.Range("E9").Value = "  -4.99%  "
$ws.Range("D10").Value = "1.099"
$ws.Range("E10").Value = "  -5.17%  "
$ws.Range("D11").Value = "0.07142"
$ws.Range("E11").Value = "  -3.94%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "21.71"
$ws.Range("E13").Value = "  -7.73%  "
$ws.Range("D14").Value = "6.059"
$ws.Range("E14").Value = "  -5.95%  "
$ws.Range("D15").Value = "1.741.76"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "6.932"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "0.06585"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "79.04"
$ws.Range("E19").Value = "  -4.58%  "
$ws.Range("D20").Value = "0.9995"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "16.69"
$ws.Range("E21").Value = "  -4.76%  "
$ws.Range("D22").Value = "6.145"
$ws.Range("E22").Value = "  -4.81%  "
$ws.Range("D23").Value = "27.757.18"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("E24").Value = "  -5.61%  "
$ws.Range("D25").Value = "2.378"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "154.03"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "19.69"
$ws.Range("E27").Value = "  -5.73%  "
$ws.Range("E28").Value = "  -6.22%  "
$ws.Range("D29").Value = "1.939.61"
$ws.Range("E30").Value = "  -12.24%  "
$ws.Range("D31").Value = "127.41"
$ws.Range("E31").Value = "  -6.51%  "
$ws.Range("D32").Value = "4.066"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").Value = "5.720"
$ws.Range("E33").Value = "  -7.31%  "
$ws.Range("D34").Value = "0.08703"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").Value = "11.88"
$ws.Range("E35").Value = "  -7.51%  "
$ws.Range("D36").Value = "1.510"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "0.02256"
$ws.Range("E37").Value = "  -7.57%  "
$ws.Range("D38").Value = "5.069"
$ws.Range("E38").Value = "  -5.29%  "
$ws.Range("D39").Value = "0.06052"
$ws.Range("E39").Value = "  -5.09%  "
$ws.Range("D40").Value = "0.6385"
$ws.Range("E40").Value = "  -7.27%  "
$ws.Range("D41").Value = "0.2077"
$ws.Range("E41").Value = "  -4.99%  "
$ws.Range("D42").Value = "1.184"
$ws.Range("E42").Value = "  -4.61%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "7.813"
$ws.Range("E44").Value = "  -6.52%  "
$ws.Range("D45").Value = "13.50"
$ws.Range("E45").Value = "  -5.76%  "
$ws.Range("D46").Value = "3.805"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").Value = "0.5893"
$ws.Range("E47").Value = "  -6.77%  "
$ws.Range("D48").Value = "125.41"
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("E49").Value = "  -6.59%  "
$ws.Range("D50").Value = "0.06926"
$ws.Range("E50").Value = "  -7.30%  "
$ws.Range("D51").Value = "1.137"
$ws.Range("E51").Value = "  -5.83%  "

# Restore the default style so the text-format override does not
# linger as a visible per-cell style change.
$ws.Range("D2:D51").Style = "Normal"
